$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 424 (existing rows 424.. shift down to 426..)
$ws.Range("A424:A425").EntireRow.Insert()

# Row 424: KODEX 차이나2차전지MSCI(합성)
$ws.Range("A424").Value = 419430
$ws.Range("B424").Value = "KODEX 차이나2차전지MSCI(합성)"
$ws.Range("C424").Value = "해외"
$ws.Range("D424").Value = "중국"

# Row 425: KODEX 미국클린에너지나스닥
$ws.Range("A425").Value = 419420
$ws.Range("B425").Value = "KODEX 미국클린에너지나스닥"
$ws.Range("C425").Value = "해외"
$ws.Range("D425").Value = "미국"

# Update the _FilterDatabase defined name range to include the 2 new rows
$fd = $wb.Names.Item(1)
$fd.RefersTo = '=Sheet1!$A$1:$D$527'

# Update view state: selection moves to B426, frozen-pane scroll follows
$ws.Activate()
$ws.Range("B426").Select()
$excel.ActiveWindow.ScrollRow = 407
